$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells receiving a numeric-looking price string must be forced to Text format
# so Excel stores them verbatim (matching the source feed's string formatting)
# instead of silently parsing them into floating point numbers.

$textCells = "D4,D5,D6,D7,D9,D10,D12,D14,D15,D17,D19,D20,D21,D23,D24,D25,D26,D27,D29,D30,D31,D32,D33,D34,D36,D37,D39,D40,D41,D42,D43,D44,D45,D47,D48,D50,D51"
foreach ($addr in $textCells.Split(",")) {
    $ws.Range($addr).NumberFormat = "@"
}

# Rows 15/16, 29/30, 41/42, 47/48: coin name + link swapped with neighboring row
$ws.Range("B15").Value = "Polkadot"
$ws.Range("C15").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D15").Value = "7.35"
$ws.Range("E15").Value = "  -0.35%  "

$ws.Range("B16").Value = "WrappedEther"
$ws.Range("C16").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D16").Value = "2.945.61"
$ws.Range("E16").Value = "  +0.37%  "

$ws.Range("B29").Value = "LEO"
$ws.Range("C29").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D29").Value = "4.10"
$ws.Range("E29").Value = "  -0.64%  "

$ws.Range("B30").Value = "Dai"
$ws.Range("C30").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D30").Value = "1.00"
$ws.Range("E30").Value = "  +0.05%  "

$ws.Range("B41").Value = "Stacks"
$ws.Range("C41").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D41").Value = "2.52"
$ws.Range("E41").Value = "  -3.09%  "

$ws.Range("B42").Value = "Stellar"
$ws.Range("C42").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D42").Value = "0.114"
$ws.Range("E42").Value = "  +0.07%  "

$ws.Range("B47").Value = "TheGraph"
$ws.Range("C47").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D47").Value = "0.270"
$ws.Range("E47").Value = "  -1.56%  "

$ws.Range("B48").Value = "ApeXProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D48").Value = "2.36"
$ws.Range("E48").Value = "  +2.03%  "

# Remaining rows: Price (D) / Volume(1h) (E) refresh only
$ws.Range("D2").Value = "51.130.48"
$ws.Range("E2").Value = "  -0.78%  "
$ws.Range("D3").Value = "2.957.31"
$ws.Range("E3").Value = "  +0.99%  "
$ws.Range("D4").Value = "0.998"
$ws.Range("E4").Value = "  -0.14%  "
$ws.Range("D5").Value = "378.54"
$ws.Range("E5").Value = "  +0.59%  "
$ws.Range("D6").Value = "102.47"
$ws.Range("E6").Value = "  -1.69%  "
$ws.Range("D7").Value = "0.538"
$ws.Range("E7").Value = "  -1.00%  "
$ws.Range("E8").Value = "  -0.02%  "
$ws.Range("D9").Value = "0.587"
$ws.Range("E9").Value = "  -0.12%  "
$ws.Range("D10").Value = "36.50"
$ws.Range("E10").Value = "  -1.24%  "
$ws.Range("E11").Value = "  -0.30%  "
$ws.Range("D12").Value = "0.0838"
$ws.Range("E12").Value = "  +0.15%  "
$ws.Range("D13").Value = "3.418.55"
$ws.Range("E13").Value = "  +0.56%  "
$ws.Range("D14").Value = "17.96"
$ws.Range("E14").Value = "  -2.23%  "
$ws.Range("D17").Value = "0.980"
$ws.Range("E17").Value = "  +4.27%  "
$ws.Range("D18").Value = "51.135.47"
$ws.Range("E18").Value = "  -0.71%  "
$ws.Range("D19").Value = "3.21"
$ws.Range("E19").Value = "  -6.14%  "
$ws.Range("D20").Value = "7.16"
$ws.Range("E20").Value = "  -2.47%  "
$ws.Range("D21").Value = "12.53"
$ws.Range("E21").Value = "  -3.59%  "
$ws.Range("D22").Value = "0.0₃0950"
$ws.Range("E22").Value = "  +0.39%  "
$ws.Range("D23").Value = "68.28"
$ws.Range("E23").Value = "  -0.20%  "
$ws.Range("D24").Value = "261.32"
$ws.Range("E24").Value = "  -0.30%  "
$ws.Range("D25").Value = "2.84"
$ws.Range("E25").Value = "  +2.21%  "
$ws.Range("D26").Value = "8.20"
$ws.Range("E26").Value = "  +11.56%  "
$ws.Range("D27").Value = "7.62"
$ws.Range("E27").Value = "  +8.37%  "
$ws.Range("E28").Value = "  -1.06%  "
$ws.Range("D31").Value = "25.66"
$ws.Range("E31").Value = "  -0.59%  "
$ws.Range("D32").Value = "0.111"
$ws.Range("E32").Value = "  +9.79%  "
$ws.Range("D33").Value = "9.77"
$ws.Range("E33").Value = "  -0.57%  "
$ws.Range("D34").Value = "50.51"
$ws.Range("E34").Value = "  -2.48%  "
$ws.Range("E35").Value = "  -2.91%  "
$ws.Range("D36").Value = "33.60"
$ws.Range("E36").Value = "  -1.61%  "
$ws.Range("D37").Value = "0.0441"
$ws.Range("E37").Value = "  +3.37%  "
$ws.Range("E38").Value = "  -0.01%  "
$ws.Range("D39").Value = "2.96"
$ws.Range("E39").Value = "  -1.78%  "
$ws.Range("D40").Value = "16.83"
$ws.Range("E40").Value = "  -0.86%  "
$ws.Range("D43").Value = "1.77"
$ws.Range("E43").Value = "  -3.10%  "
$ws.Range("D44").Value = "121.91"
$ws.Range("E44").Value = "  -2.05%  "
$ws.Range("D45").Value = "20.97"
$ws.Range("E45").Value = "  -4.06%  "
$ws.Range("E46").Value = "  -0.75%  "
$ws.Range("D49").Value = "2.002.43"
$ws.Range("E49").Value = "  -1.06%  "
$ws.Range("D50").Value = "3.20"
$ws.Range("E50").Value = "  +0.68%  "
$ws.Range("D51").Value = "0.0332"
$ws.Range("E51").Value = "  +3.50%  "
